$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down.
$ws.Rows.Item(1).Insert()

# Add header labels (stored as shared strings) in the new row 1.
$ws.Range("A1").Value = "Cutting speed V"
$ws.Range("B1").Value = "Feed per tooth"
$ws.Range("C1").Value = "Axial depth"
$ws.Range("E1").Value = "tool life"

# Column widths tweaked by the author (closest reachable values; the
# engine snaps ColumnWidth to whole-pixel increments).
$ws.Columns.Item(1).ColumnWidth = 13.29
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 11

# Update the active selection.
$ws.Range("I3").Select()
